# After completion of SearchAndSaveCastTest
# - Input sheet: move the remembered selection from B3 to B1
# - Series Cast sheet: append the two new cast members captured by the test
#   (Toby Foster/Dennis and Julia Dearden/Nun), then leave the sheet
#   scrolled back to the top with A3:C65 selected (matching the new data
#   range) as the active sheet.

$wb = $excel.ActiveWorkbook

# --- "Input" sheet: just update the remembered selection -------------------
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("B1").Select()

# --- "Series Cast" sheet: append newly scraped rows -------------------------
$wsCast = $wb.Worksheets.Item("Series Cast")

$wsCast.Range("A50").Value = "Toby Foster"
$wsCast.Range("B50").Value = "Dennis"
$wsCast.Range("C50").Value = "1 episode, 2020"

$wsCast.Range("A51").Value = "Julia Dearden"
$wsCast.Range("B51").Value = "Nun"
$wsCast.Range("C51").Value = "1 episode, 2019"

# Leave this sheet active/selected, as it was before the edit, with the
# updated selection covering the (now larger) data range.
$wsCast.Range("A3:C65").Select()
